# Rename the first worksheet ("Monkey" -> "This be a sheet") and make sure
# it is the active/selected sheet (the workbook's active tab moves from
# "SheetA" back onto the renamed first sheet).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "This be a sheet"
$ws1.Activate()

# SheetA's page setup no longer forces a custom first page number.
$ws2 = $wb.Worksheets.Item(2)
$ws2.PageSetup.FirstPageNumber = 0
